$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2" = "261.12"
    "D3" = "26.97"
    "E3" = "-0.74%"
    "D4" = "4.716"
    "E4" = "-0.27%"
    "E5" = "2.24%"
    "D6" = "6.751"
    "E6" = "1.26%"
    "D7" = "0.8514"
    "E7" = "0.68%"
    "D8" = "0.9147"
    "E8" = "-0.79%"
    "D9" = "0.1402"
    "E9" = "-0.44%"
    "D10" = "0.05082"
    "E10" = "3.00%"
    "D11" = "0.07066"
    "E11" = "-0.40%"
    "D12" = "0.03098"
    "E12" = "-1.22%"
    "D13" = "0.09052"
    "E13" = "-0.27%"
    "D14" = "0.001540"
    "E14" = "-0.38%"
    "D15" = "0.0006172"
    "E15" = "1.00%"
    "D16" = "0.005981"
    "E16" = "-2.41%"
    "D17" = "3.445"
    "E17" = "-0.22%"
    "D18" = "3.167"
    "E18" = "0.38%"
    "E19" = "-1.04%"
    "E20" = "-0.68%"
    "D21" = "0.1310"
    "E21" = "1.03%"
    "D22" = "4.080"
    "E22" = "-0.25%"
    "D23" = "0.04250"
    "E23" = "0.04%"
    "D24" = "0.001200"
    "E24" = "-1.57%"
    "D25" = "0.004076"
    "E25" = "4.20%"
    "D26" = "0.0001200"
    "E26" = "0.04%"
    "E27" = "4.15%"
    "D40" = "0.03956"
    "E40" = "2.04%"
    "E41" = "-0.06%"
    "D42" = "0.004139"
    "D44" = "0.01358"
    "E44" = "-16.90%"
    "E45" = "-3.07%"
    "E46" = "0.04%"
    "D48" = "0.2531"
    "E48" = "91.61%"
    "E49" = "0.04%"
    "E50" = "0.04%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
